$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source inlineStr cells)
$ws.Range('D2').Value = '60.388.25'
$ws.Range('E2').Value = '  -4.21%  '
$ws.Range('D3').Value = '2.992.56'
$ws.Range('E3').Value = '  -5.51%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.17'
$ws.Range('E5').Value = '  -3.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '125.56'
$ws.Range('E6').Value = '  -6.73%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').Value = '2.985.37'
$ws.Range('E8').Value = '  -5.62%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.498'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('E10').Value = '  -5.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.09'
$ws.Range('E11').Value = '  -3.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.436'
$ws.Range('E12').Value = '  -3.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -5.68%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.57'
$ws.Range('E14').Value = '  -6.50%  '
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').Value = '3.481.82'
$ws.Range('E16').Value = '  -5.62%  '
$ws.Range('D17').Value = '60.335.52'
$ws.Range('E17').Value = '  -4.23%  '
$ws.Range('D18').Value = '2.982.48'
$ws.Range('E18').Value = '  -5.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.15'
$ws.Range('E19').Value = '  -6.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.52'
$ws.Range('E20').Value = '  -6.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.05'
$ws.Range('E21').Value = '  -6.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.659'
$ws.Range('E22').Value = '  -5.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').Value = '  -6.52%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.94'
$ws.Range('E24').Value = '  -3.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.70'
$ws.Range('E25').Value = '  -5.59%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.50'
$ws.Range('E28').Value = '  -6.46%  '
$ws.Range('E29').Value = '  -7.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.89'
$ws.Range('E30').Value = '  -6.86%  '
$ws.Range('E31').Value = '  -7.11%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.01'
$ws.Range('E32').Value = '  -10.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0932'
$ws.Range('E33').Value = '  -9.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.26'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('E35').Value = '  -8.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.54'
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '49.34'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('D38').Value = '0.0₃0650'
$ws.Range('E38').Value = '  -8.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0357'
$ws.Range('E39').Value = '  -8.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.82'
$ws.Range('E40').Value = '  -4.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.106'
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '373.28'
$ws.Range('E42').Value = '  -7.85%  '
$ws.Range('D43').Value = '2.661.27'
$ws.Range('E43').Value = '  -5.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.39'
$ws.Range('E44').Value = '  -6.91%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('E46').Value = '  -7.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '118.07'
$ws.Range('E47').Value = '  -4.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.95'
$ws.Range('E48').Value = '  -8.22%  '
$ws.Range('E49').Value = '  -4.88%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.35'
$ws.Range('E50').Value = '  -7.88%  '
$ws.Range('B51').Value = 'Arweave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.96'
$ws.Range('E51').Value = '  -5.61%  '
